$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A width (target stored width 14.42578125 chars; the
# ColumnWidth setter snaps to a pixel grid, so feed it the value that
# rounds to the closest achievable stored width)
$ws.Columns.Item(1).ColumnWidth = 13.666666666666666

# Update cell values
$ws.Range("A1").Value = -0.043659189892190989
$ws.Range("B1").Value = -0.042454134145798315

$ws.Range("A2").Value = -0.016961780446142133
$ws.Range("B2").Value = -0.012601184040415367

$ws.Range("A3").Value = -0.018181246521458031
$ws.Range("B3").Value = -0.0060624633259771391

$ws.Range("A4").Value = -0.077187587239905278
$ws.Range("B4").Value = -0.076486373405838345

$ws.Range("A5").Value = -0.034728083286631205
$ws.Range("B5").Value = -0.014122491216763487
